# Apply the "Chuan bi cham do an" update to the employee roster.
# - Row for employee #3 (row 4) gets new personal details and a new sales role.
# - Several "Nhân viên quản lý kho" roles are re-classified into more specific
#   roles (sales staff, position manager, store manager, customer care), each
#   with an updated Final Salary figure.
# - Employee #12 (row 13) is reactivated ("Hoạt động").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (ID = 3): new hire details replacing the old placeholder employee.
$ws.Cells.Item(4, 2).Value = "Tần Thiên"
$ws.Cells.Item(4, 3).Value = "Lang"
$ws.Cells.Item(4, 4).Value = "15-01-2000"
$ws.Cells.Item(4, 5).Value = "Nhân viên bán hàng"
$ws.Cells.Item(4, 7).Value = "6.000.000"

# Row 5 (ID = 4): re-classified as sales staff, final salary recomputed.
$ws.Cells.Item(5, 5).Value = "Nhân viên bán hàng"
$ws.Cells.Item(5, 7).Value = "7.000.000"

# Row 6 (ID = 5): re-classified as sales staff, final salary recomputed.
$ws.Cells.Item(6, 5).Value = "Nhân viên bán hàng"
$ws.Cells.Item(6, 7).Value = "8.000.000"

# Row 7 (ID = 6): re-classified as sales staff, final salary recomputed.
$ws.Cells.Item(7, 5).Value = "Nhân viên bán hàng"
$ws.Cells.Item(7, 7).Value = "9.000.000"

# Row 8 (ID = 7): re-classified as position manager, final salary recomputed.
$ws.Cells.Item(8, 5).Value = "Quản lý chức vụ"
$ws.Cells.Item(8, 7).Value = "7.840.000"

# Row 9 (ID = 8): re-classified as store manager, final salary recomputed.
$ws.Cells.Item(9, 5).Value = "Quản lý cửa hàng"
$ws.Cells.Item(9, 7).Value = "8.320.000"

# Row 10 (ID = 9): re-classified as customer care staff, final salary recomputed.
$ws.Cells.Item(10, 5).Value = "Nhân viên chăm sóc khách hàng"
$ws.Cells.Item(10, 7).Value = "6.380.000"

# Row 13 (ID = 12): status flipped back to active.
$ws.Cells.Item(13, 8).Value = "Hoạt động"

# The Role Name column now contains a longer string ("Nhân viên chăm sóc
# khách hàng"), so re-fit the column width to the new widest entry.
$ws.Columns.Item(5).ColumnWidth = 28.75
